# meetings.xlsx: add a "Comment (Optional)" column (E) with a sample row,
# and correct the example Zoom link's displayed text.
#
# NOTE on ordering: new shared-strings are interned in the order values are
# first written, and the target file expects "Example meeting from Excel"
# before "Comment (Optional)" in xl/sharedStrings.xml, so E2 is set before E1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample row value for the new column (written first so it gets the
# earlier shared-string index, matching the authored file).
$ws.Range("E2").Value = "Example meeting from Excel"

# New column header.
$ws.Range("E1").Value = "Comment (Optional)"

# Fix up the displayed Zoom link text (the hyperlink target itself is left
# untouched, matching the upstream edit).
$ws.Range("B2").Value = "https://us05web.zoom.us/j/87177504375?pwd=jhvL2kxa2ZoQWdicWd1BiS0JLZzadblahblah"

# Move the active selection to B2 (was B9 before the edit).
$ws.Range("B2").Select()
